# Updates cryptos list cell values per the authoritative diff.
# Plain text / URL / percentage-string cells: assigned directly.
# Price cells whose new text is a bare number (e.g. "1.00", "6.22") would be
# auto-coerced to a numeric Value by Excel's normal typed-input parsing, which
# would silently lose the literal formatting ("1.00" -> 1). Those cells are
# marked as Text (NumberFormat "@") immediately before the write so the exact
# source string is preserved, matching the original inline-string cell contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textUpdates = [ordered]@{
    'D2' = '67.878.18'
    'E2' = '  -0.84%  '
    'D3' = '3.254.13'
    'E3' = '  -0.61%  '
    'E4' = '  -0.02%  '
    'E5' = '  -0.16%  '
    'E6' = '  -0.64%  '
    'E7' = '  +0.08%  '
    'E8' = '  +0.08%  '
    'E9' = '  -3.46%  '
    'E10' = '  -0.81%  '
    'E11' = '  -3.88%  '
    'D12' = '3.825.84'
    'E12' = '  -0.41%  '
    'E13' = '  +1.60%  '
    'B14' = 'WrappedBTC'
    'C14' = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
    'D14' = '67.923.70'
    'E14' = '  -0.79%  '
    'B15' = 'Avalanche'
    'C15' = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
    'E15' = '  -4.27%  '
    'E16' = '  -2.73%  '
    'D17' = '3.313.15'
    'E17' = '  +1.16%  '
    'E18' = '  -2.63%  '
    'E19' = '  -2.94%  '
    'E20' = '  +5.85%  '
    'E21' = '  -3.20%  '
    'E22' = '  -0.14%  '
    'E23' = '  -1.01%  '
    'E24' = '  -2.90%  '
    'E25' = '  -3.21%  '
    'E26' = '  -1.05%  '
    'E27' = '  -3.30%  '
    'E28' = '  +0.26%  '
    'E29' = '  -2.05%  '
    'E30' = '  -2.43%  '
    'E31' = '  -5.07%  '
    'E32' = '  -5.29%  '
    'E33' = '  -5.00%  '
    'E34' = '  -0.83%  '
    'E35' = '  -5.53%  '
    'E36' = '  -3.36%  '
    'E37' = '  +0.14%  '
    'E38' = '  -4.47%  '
    'E39' = '  -4.11%  '
    'E40' = '  -5.21%  '
    'D41' = '2.631.96'
    'E41' = '  -0.74%  '
    'B42' = 'Hedera'
    'C42' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'E42' = '  -2.62%  '
    'B43' = 'dogwifhat'
    'C43' = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
    'E43' = '  -5.70%  '
    'B44' = 'Bittensor'
    'C44' = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
    'E44' = '  -1.76%  '
    'B45' = 'InjectiveProtocol'
    'C45' = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
    'E45' = '  -4.61%  '
    'B46' = 'VeChain'
    'C46' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'E46' = '  -3.67%  '
    'B47' = 'Cosmos'
    'C47' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'E47' = '  -2.25%  '
    'B48' = 'ONDO'
    'C48' = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
    'E48' = '  -2.26%  '
    'B49' = 'Stellar'
    'C49' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'E49' = '  -2.37%  '
    'B50' = 'FirstDigitalUSD'
    'C50' = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
    'E50' = '  -0.12%  '
    'B51' = 'Arweave'
    'C51' = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
    'E51' = '  -4.96%  '
}

$numericTextUpdates = [ordered]@{
    'D5' = '581.39'
    'D6' = '183.97'
    'D8' = '0.599'
    'D10' = '6.60'
    'D11' = '0.406'
    'D15' = '27.23'
    'D18' = '5.69'
    'D19' = '13.21'
    'D20' = '414.62'
    'D21' = '7.50'
    'D22' = '1.00'
    'D23' = '71.12'
    'D24' = '0.506'
    'D25' = '0.0000116'
    'D26' = '0.186'
    'D27' = '9.36'
    'D30' = '22.57'
    'D31' = '5.43'
    'D32' = '6.82'
    'D34' = '162.89'
    'D35' = '1.43'
    'D36' = '1.87'
    'D37' = '26.89'
    'D38' = '0.792'
    'D39' = '4.43'
    'D40' = '6.28'
    'D42' = '0.0672'
    'D43' = '2.41'
    'D44' = '336.84'
    'D45' = '24.13'
    'D46' = '0.0272'
    'D47' = '6.22'
    'D48' = '0.973'
    'D49' = '0.0999'
    'D50' = '1.00'
    'D51' = '30.43'
}

foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

foreach ($ref in $numericTextUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = '@'
    $cell.Value = $numericTextUpdates[$ref]
}

